$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows 45-54 appended to the "20000-40000" TCV range sheet.
# Column order: A=CVR, B=Year, C=Beløb 12 mdr. (TCV), D=Løsning,
#               E=Opsagt dato:, F=Årsag (unused here), G=Ny leverandør,
#               H=Quarter, I=TCV_range
$newRows = @(
    @{ Row=45; CVR='21699780'; Year=2023; TCV=21650.08; Solution='Visma Løn'; Date=45182; NewVendor='ADP Celergo';              Quarter='2023Q3' },
    @{ Row=46; CVR='40881239'; Year=2023; TCV=28521.38; Solution='Visma Løn'; Date=45205; NewVendor='DataLøn';                  Quarter='2023Q4' },
    @{ Row=47; CVR='67437853'; Year=2023; TCV=30097;    Solution='Visma Løn'; Date=45205; NewVendor='DataLøn';                  Quarter='2023Q4' },
    @{ Row=48; CVR='33738811'; Year=2023; TCV=39408;    Solution='Visma Time'; Date=45211; NewVendor=$null;                     Quarter='2023Q4' },
    @{ Row=49; CVR='28155379'; Year=2023; TCV=20165;    Solution='SKL SE';     Date=45239; NewVendor=$null;                     Quarter='2023Q4' },
    @{ Row=50; CVR='29186898'; Year=2023; TCV=23777;    Solution='Visma Løn'; Date=45257; NewVendor=$null;                     Quarter='2023Q4' },
    @{ Row=51; CVR='47252059'; Year=2023; TCV=31857;    Solution='Visma Løn'; Date=45245; NewVendor=$null;                     Quarter='2023Q4' },
    @{ Row=52; CVR='14364641'; Year=2023; TCV=22553;    Solution='Visma Løn'; Date=45268; NewVendor='DataLøn';                  Quarter='2023Q4' },
    @{ Row=53; CVR='19764338'; Year=2023; TCV=24550;    Solution='SKL SE';     Date=45281; NewVendor=$null;                     Quarter='2023Q4' },
    @{ Row=54; CVR='28979355'; Year=2023; TCV=20119;    Solution='Visma Løn'; Date=45289; NewVendor=$null;                     Quarter='2023Q4' }
)

foreach ($r in $newRows) {
    $row = $r.Row

    # Column A (CVR) - looks numeric, force text storage like the other CVR
    # cells in the sheet (quote-prefix trick), then drop the quote-prefix
    # style so the cell ends up styleless like its siblings.
    $cellA = $ws.Cells.Item($row, 1)
    $cellA.Value = "'" + $r.CVR
    $cellA.Style = "Normal"

    # Column B (Year)
    $ws.Cells.Item($row, 2).Value = $r.Year

    # Column C (Beløb 12 mdr. (TCV))
    $ws.Cells.Item($row, 3).Value = $r.TCV

    # Column D (Løsning)
    $ws.Cells.Item($row, 4).Value = $r.Solution

    # Column E (Opsagt dato:) - stored as a date serial with the sheet's
    # existing custom date/time number format.
    $cellE = $ws.Cells.Item($row, 5)
    $cellE.Value = $r.Date
    $cellE.NumberFormat = "YYYY-MM-DD HH:MM:SS"

    # Column G (Ny leverandør) - only populated for some rows
    if ($r.NewVendor) {
        $ws.Cells.Item($row, 7).Value = $r.NewVendor
    }

    # Column H (Quarter)
    $ws.Cells.Item($row, 8).Value = $r.Quarter

    # Column I (TCV_range)
    $ws.Cells.Item($row, 9).Value = "20000-40000"
}
